{"js": "// Replace the date line and all 25 multiplication problems with their\n// new values, as described by the diff. Every old text value in this\n// document is unique, so a simple search + replace per pair is safe.\nconst replacements = [\n  [\"2023-09-16 Saturday\", \"2023-09-17 Sunday\"],\n  [\"52\u00d721=\", \"76\u00d761=\"],\n  [\"68\u00d796=\", \"39\u00d763=\"],\n  [\"31\u00d771=\", \"80\u00d768=\"],\n  [\"76\u00d740=\", \"44\u00d728=\"],\n  [\"83\u00d762=\", \"65\u00d754=\"],\n  [\"67\u00d788=\", \"50\u00d765=\"],\n  [\"61\u00d767=\", \"56\u00d777=\"],\n  [\"40\u00d721=\", \"30\u00d743=\"],\n  [\"17\u00d793=\", \"86\u00d762=\"],\n  [\"15\u00d794=\", \"97\u00d783=\"],\n  [\"72\u00d780=\", \"74\u00d790=\"],\n  [\"57\u00d753=\", \"53\u00d751=\"],\n  [\"28\u00d738=\", \"63\u00d764=\"],\n  [\"75\u00d775=\", \"93\u00d767=\"],\n  [\"62\u00d717=\", \"75\u00d773=\"],\n  [\"63\u00d782=\", \"15\u00d776=\"],\n  [\"23\u00d765=\", \"40\u00d724=\"],\n  [\"90\u00d793=\", \"99\u00d760=\"],\n  [\"84\u00d721=\", \"53\u00d719=\"],\n  [\"11\u00d717=\", \"27\u00d795=\"],\n  [\"90\u00d720=\", \"56\u00d731=\"],\n  [\"93\u00d729=\", \"16\u00d734=\"],\n  [\"75\u00d795=\", \"45\u00d717=\"],\n  [\"94\u00d735=\", \"85\u00d733=\"],\n  [\"48\u00d732=\", \"70\u00d737=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all 25 multiplication problems with their\n# new values, as described by the diff. Every old text value in this\n# document is unique, so a simple Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2023-09-16 Saturday\", \"2023-09-17 Sunday\")\n  ,@(\"52\u00d721=\", \"76\u00d761=\")\n  ,@(\"68\u00d796=\", \"39\u00d763=\")\n  ,@(\"31\u00d771=\", \"80\u00d768=\")\n  ,@(\"76\u00d740=\", \"44\u00d728=\")\n  ,@(\"83\u00d762=\", \"65\u00d754=\")\n  ,@(\"67\u00d788=\", \"50\u00d765=\")\n  ,@(\"61\u00d767=\", \"56\u00d777=\")\n  ,@(\"40\u00d721=\", \"30\u00d743=\")\n  ,@(\"17\u00d793=\", \"86\u00d762=\")\n  ,@(\"15\u00d794=\", \"97\u00d783=\")\n  ,@(\"72\u00d780=\", \"74\u00d790=\")\n  ,@(\"57\u00d753=\", \"53\u00d751=\")\n  ,@(\"28\u00d738=\", \"63\u00d764=\")\n  ,@(\"75\u00d775=\", \"93\u00d767=\")\n  ,@(\"62\u00d717=\", \"75\u00d773=\")\n  ,@(\"63\u00d782=\", \"15\u00d776=\")\n  ,@(\"23\u00d765=\", \"40\u00d724=\")\n  ,@(\"90\u00d793=\", \"99\u00d760=\")\n  ,@(\"84\u00d721=\", \"53\u00d719=\")\n  ,@(\"11\u00d717=\", \"27\u00d795=\")\n  ,@(\"90\u00d720=\", \"56\u00d731=\")\n  ,@(\"93\u00d729=\", \"16\u00d734=\")\n  ,@(\"75\u00d795=\", \"45\u00d717=\")\n  ,@(\"94\u00d735=\", \"85\u00d733=\")\n  ,@(\"48\u00d732=\", \"70\u00d737=\")\n)\n\nforeach ($pair in $pairs) {\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
